$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 2 (the "2007年" data row). This shifts rows 3-6
# up to become rows 2-5, matching the diff (2010/2012/2015/2017年 rows
# moving up one position) and shrinking the used range from A1:R6 to A1:R5.
$ws.Rows.Item(2).Delete()
